$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The CasesTab query (B2) had its trailing "Cohort" output column removed
# (the OPTIONAL MATCH (co:cohort) line stays, only the returned coalesce(...) AS `Cohort` column is dropped).
$ws.Range("B2").Value = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`r`n`r`nMATCH (c)<--(diag:diagnosis)`r`nWHERE s.clinical_study_designation IN ['UBC01'] and demo.breed in ['Border Collie','Chihuahua','Maltese','Scottish Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in ['Bladder, Urethra', 'Bladder, Urethra, Prostate']`r`nOPTIONAL MATCH (samp:sample)-->(c)`r`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`r`nWITH DISTINCT c, s, demo, diag, co`r`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`r`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`r`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`r`n        coalesce(demo.breed, '') AS Breed ,`r`n        coalesce(diag.disease_term, '') AS Diagnosis ,`r`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`r`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`r`n        coalesce(demo.sex, '') AS Sex ,`r`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`r`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`r`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

# Shorter text means the wrapped row shrinks a bit.
$ws.Rows.Item(2).RowHeight = 304.5

# The selection moved from the old B8/row4 scroll position back up to B2.
$ws.Range("B2").Select()
